# Apply the "front end paginas alertas, mensagens, tarefas e recuperação de senha"
# update to the Cronograma worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

# D15 - mark "Criar tela de alteração de senha" task as "Stand by"
$ws.Range("D15").Value = "Stand by"

# New rows 21-23: front end pages for alerts, messages and tasks
$ws.Range("B21").Value = "Criar tela de alertas"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = "ok"

$ws.Range("B22").Value = "Criar tela de mensagens"
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = "ok"

$ws.Range("B23").Value = "Criar tela de tarefas"
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = "ok"

# Update the view: select B23 (matches the new last filled row)
$ws.Activate()
$ws.Range("B23").Select()
